$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value  = 3.8
$ws.Range("H3").Value  = 3.5
$ws.Range("J3").Value  = 4.3
$ws.Range("L3").Value  = 2.47
$ws.Range("M3").Value  = 1.06
$ws.Range("N3").Value  = 7.5
$ws.Range("O3").Value  = 1.28
$ws.Range("P3").Value  = 3.35
$ws.Range("Q3").Value  = 1.83
$ws.Range("R3").Value  = 1.91
$ws.Range("S3").Value  = 1.42
$ws.Range("T3").Value  = 2.67
$ws.Range("V3").Value  = 2
$ws.Range("W3").Value  = 11.25
$ws.Range("AA3").Value = 35
$ws.Range("AB3").Value = 40
$ws.Range("AC3").Value = 7.5
$ws.Range("AO3").Value = 22
$ws.Range("AP3").Value = 29
$ws.Range("AQ3").Value = 120
$ws.Range("AR3").Value = 175
$ws.Range("AS3").Value = 400
$ws.Range("AT3").Value = 2.67
$ws.Range("AU3").Value = 7.3
$ws.Range("AV3").Value = 70
$ws.Range("AX3").Value = 3.75
